# 孫大千 財產申報表 — "#5: cash & deposit done"
# Adds bank/deposit_type/currency metadata columns (property_category, category,
# date, legislator_name, legislator_id, source_file, index) to the 存款 (deposit)
# sheet, matching the layout already used on the other sheets (股票, 具有相當價值之財產, 保險).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# --- Header row (row 1) -----------------------------------------------------
# B1:D1 already hold "bank"/"deposit_type"/"currency" labels; E1/F1 need to move
# from literal sample values to the "owner"/"total" labels, and new G1:M1 labels
# are appended for the metadata columns.
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# Match the existing bold/bordered header styling on the newly added header cells.
$ws.Range("B1").Copy()
$ws.Range("G1:M1").PasteSpecial(-4122)

# --- Data rows (rows 2-11) ---------------------------------------------------
# Columns A-F are untouched; append the same constant metadata block used by the
# other property sheets, with the running index mirrored into the new column M.
for ($r = 2; $r -le 11; $r++) {
    $idx = $ws.Range("A$r").Value()

    $ws.Range("G$r").Value = "deposit"
    $ws.Range("H$r").Value = "normal"
    $ws.Range("I$r").Value = "2012-03-03"
    $ws.Range("J$r").Value = "孫大千"
    $ws.Range("K$r").Value = 919
    $ws.Range("L$r").Value = "tmpc261"
    $ws.Range("M$r").Value = $idx
}
